# Auto-generated: refresh market-price-derived columns (H-N) on the 8 job-leve
# profit sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR, CUL) to match a fresh pull
# of Universalis current-average-price data. Values only; no formulas, styles,
# rows, or table definitions change.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H33").Value = 198.2963
$ws.Range("I33").Value = 114.10526
$ws.Range("K33").Value = 114.10526
$ws.Range("M33").Value = 114.89474

$ws.Range("H64").Value = 28022.445
$ws.Range("I64").Value = 34071.715
$ws.Range("K64").Value = 34071.715
$ws.Range("M64").Value = -33823.715

$ws.Range("H67").Value = 28022.445
$ws.Range("I67").Value = 34071.715
$ws.Range("K67").Value = 34071.715
$ws.Range("M67").Value = -33213.715

$ws.Range("H88").Value = 6189.4546
$ws.Range("J88").Value = 10738.8
$ws.Range("L88").Value = 10738.8
$ws.Range("N88").Value = -11550.8

$ws.Range("H91").Value = 6189.4546
$ws.Range("J91").Value = 10738.8
$ws.Range("L91").Value = 10738.8
$ws.Range("N91").Value = -13546.8

$ws.Range("H92").Value = 1661.36
$ws.Range("I92").Value = 167.05263
$ws.Range("K92").Value = 167.05263
$ws.Range("M92").Value = 1080.94737

$ws.Range("H106").Value = 1374.091
$ws.Range("I106").Value = 1374.091
$ws.Range("K106").Value = 1374.091
$ws.Range("M106").Value = -743.0909999999999

$ws.Range("H116").Value = 14150.156
$ws.Range("I116").Value = 16329.792
$ws.Range("J116").Value = 7611.25
$ws.Range("K116").Value = 16329.792
$ws.Range("L116").Value = 7611.25
$ws.Range("M116").Value = -12887.792
$ws.Range("N116").Value = -14495.25

$ws.Range("H132").Value = 24332.902
$ws.Range("I132").Value = 27474.814
$ws.Range("K132").Value = 82424.442
$ws.Range("M132").Value = -79894.442

$ws.Range("H135").Value = 3336.5334
$ws.Range("I135").Value = 3237.7827
$ws.Range("K135").Value = 29140.0443
$ws.Range("M135").Value = -26605.0443

$ws.Range("H137").Value = 26702.334
$ws.Range("I137").Value = 52332.832
$ws.Range("J137").Value = 1071.8334
$ws.Range("K137").Value = 156998.496
$ws.Range("L137").Value = 3215.5002
$ws.Range("M137").Value = -154448.496
$ws.Range("N137").Value = -8315.5002

$ws.Range("H141").Value = 3877.6667
$ws.Range("I141").Value = 3877.6667
$ws.Range("K141").Value = 11633.0001
$ws.Range("M141").Value = -6453.000100000001

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 21232.06
$ws.Range("I32").Value = 22956.695
$ws.Range("K32").Value = 22956.695
$ws.Range("M32").Value = -22669.695

$ws.Range("H45").Value = 3893.7334
$ws.Range("I45").Value = 2023.7142
$ws.Range("K45").Value = 2023.7142
$ws.Range("M45").Value = -1646.7142

$ws.Range("H61").Value = 4814.241
$ws.Range("I61").Value = 1066.625
$ws.Range("K61").Value = 1066.625
$ws.Range("M61").Value = -854.625

$ws.Range("H88").Value = 3747.3333
$ws.Range("I88").Value = 636.3333
$ws.Range("K88").Value = 636.3333
$ws.Range("M88").Value = -230.3333

$ws.Range("H91").Value = 3747.3333
$ws.Range("I91").Value = 636.3333
$ws.Range("K91").Value = 636.3333
$ws.Range("M91").Value = 767.6667

$ws.Range("H97").Value = 1733.7368
$ws.Range("I97").Value = 1354.4166
$ws.Range("J97").Value = 2384
$ws.Range("K97").Value = 1354.4166
$ws.Range("L97").Value = 2384
$ws.Range("M97").Value = -858.4166
$ws.Range("N97").Value = -3376

$ws.Range("H102").Value = 1888.6072
$ws.Range("I102").Value = 1534.04
$ws.Range("J102").Value = 4843.3335
$ws.Range("K102").Value = 1534.04
$ws.Range("L102").Value = 4843.3335
$ws.Range("M102").Value = 87.96000000000004
$ws.Range("N102").Value = -8087.3335

$ws.Range("H136").Value = 4814.241
$ws.Range("I136").Value = 1066.625
$ws.Range("K136").Value = 3199.875
$ws.Range("M136").Value = -649.875

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H99").Value = 1522.125
$ws.Range("I99").Value = 1212.8334
$ws.Range("J99").Value = 2450
$ws.Range("K99").Value = 1212.8334
$ws.Range("L99").Value = 2450
$ws.Range("M99").Value = 285.1666
$ws.Range("N99").Value = -5446

$ws.Range("H105").Value = 2025.24
$ws.Range("I105").Value = 1336.8823
$ws.Range("K105").Value = 1336.8823
$ws.Range("M105").Value = 410.1177

$ws.Range("H107").Value = 2445.6333
$ws.Range("I107").Value = 2280.3
$ws.Range("K107").Value = 2280.3
$ws.Range("M107").Value = -360.3000000000002

$ws.Range("H134").Value = 2739.2354
$ws.Range("I134").Value = 2397.6428
$ws.Range("J134").Value = 4333.3335
$ws.Range("K134").Value = 7192.928400000001
$ws.Range("L134").Value = 13000.0005
$ws.Range("M134").Value = -4657.928400000001
$ws.Range("N134").Value = -18070.0005

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H69").Value = 50000
$ws.Range("I69").Value = 50000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 50000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -49251
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 50000
$ws.Range("I72").Value = 50000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 150000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -146256
$ws.Range("N72").ClearContents()

$ws.Range("H99").Value = 7025
$ws.Range("J99").Value = 9598.5
$ws.Range("L99").Value = 9598.5
$ws.Range("N99").Value = -12594.5

$ws.Range("H105").Value = 21306.578
$ws.Range("I105").Value = 26314.2
$ws.Range("K105").Value = 26314.2
$ws.Range("M105").Value = -24567.2

$ws.Range("H107").Value = 1103.3846
$ws.Range("I107").Value = 986.7143
$ws.Range("K107").Value = 986.7143
$ws.Range("M107").Value = 933.2857

$ws.Range("H126").Value = 7025
$ws.Range("J126").Value = 9598.5
$ws.Range("L126").Value = 28795.5
$ws.Range("N126").Value = -33735.5

$ws.Range("H132").Value = 54455.21
$ws.Range("J132").Value = 4366.6
$ws.Range("L132").Value = 13099.8
$ws.Range("N132").Value = -18159.8

$ws.Range("H134").Value = 3804.7896
$ws.Range("I134").Value = 3108.8333
$ws.Range("K134").Value = 9326.499899999999
$ws.Range("M134").Value = -6791.499899999999

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H132").Value = 2012.0264
$ws.Range("I132").Value = 1766.3914
$ws.Range("K132").Value = 5299.174199999999
$ws.Range("M132").Value = -2769.174199999999

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1559.8889
$ws.Range("I22").Value = 1149.4
$ws.Range("J22").Value = 2073
$ws.Range("K22").Value = 1149.4
$ws.Range("L22").Value = 2073
$ws.Range("M22").Value = -854.4000000000001
$ws.Range("N22").Value = -2663

$ws.Range("H27").Value = 1559.8889
$ws.Range("I27").Value = 1149.4
$ws.Range("J27").Value = 2073
$ws.Range("K27").Value = 1149.4
$ws.Range("L27").Value = 2073
$ws.Range("M27").Value = -1042.4
$ws.Range("N27").Value = -2287

$ws.Range("H40").Value = 3767
$ws.Range("I40").Value = 3448
$ws.Range("K40").Value = 3448
$ws.Range("M40").Value = -3312

$ws.Range("H63").Value = 112499.75
$ws.Range("J63").Value = 112499.75
$ws.Range("L63").Value = 112499.75
$ws.Range("N63").Value = -113997.75

$ws.Range("H66").Value = 112499.75
$ws.Range("J66").Value = 112499.75
$ws.Range("L66").Value = 337499.25
$ws.Range("N66").Value = -344987.25

$ws.Range("H68").Value = 4064.75
$ws.Range("I68").Value = 3115.7778
$ws.Range("J68").Value = 5284.857
$ws.Range("K68").Value = 3115.7778
$ws.Range("L68").Value = 5284.857
$ws.Range("M68").Value = -2366.7778
$ws.Range("N68").Value = -6782.857

$ws.Range("H71").Value = 4064.75
$ws.Range("I71").Value = 3115.7778
$ws.Range("J71").Value = 5284.857
$ws.Range("K71").Value = 15578.889
$ws.Range("L71").Value = 26424.285
$ws.Range("M71").Value = -11834.889
$ws.Range("N71").Value = -33912.285

$ws.Range("H74").Value = 150000
$ws.Range("J74").Value = 150000
$ws.Range("L74").Value = 150000
$ws.Range("N74").Value = -151996

$ws.Range("H77").Value = 150000
$ws.Range("J77").Value = 150000
$ws.Range("L77").Value = 450000
$ws.Range("N77").Value = -459984

$ws.Range("H93").Value = 1308.75
$ws.Range("I93").Value = 1200.4375
$ws.Range("J93").Value = 1742
$ws.Range("K93").Value = 1200.4375
$ws.Range("L93").Value = 1742
$ws.Range("M93").Value = 47.5625
$ws.Range("N93").Value = -4238

$ws.Range("H132").Value = 4709.2856
$ws.Range("I132").Value = 4709.2856
$ws.Range("K132").Value = 14127.8568
$ws.Range("M132").Value = -11597.8568

$ws.Range("H136").Value = 4103.5713
$ws.Range("I136").Value = 3592.5715
$ws.Range("K136").Value = 10777.7145
$ws.Range("M136").Value = -8227.7145

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 63310.59
$ws.Range("I122").Value = 76608.39
$ws.Range("K122").Value = 229825.17
$ws.Range("M122").Value = -227375.17

$ws.Range("H132").Value = 21148.584
$ws.Range("I132").Value = 24391.633
$ws.Range("K132").Value = 73174.899
$ws.Range("M132").Value = -70644.899

